# Updates cryptos list: prices and volume(1h) percentages for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" (D) and "Volume(1h)" (E) columns store plain text values
# (e.g. "239.78", "  +2.77%  ") even though some look numeric. Force the
# cells to text formatting first so Excel doesn't silently convert
# numeric-looking strings into real numbers, then restore the default
# "Normal" style so no stray formatting is left behind on the cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.608.48'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '1.894.82'
$ws.Range("E3").Value = '  +0.80%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '239.78'
$ws.Range("E5").Value = '  +2.77%  '
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").Value = '0.4926'
$ws.Range("E7").Value = '  +1.25%  '
$ws.Range("D8").Value = '0.2948'
$ws.Range("E8").Value = '  +2.50%  '
$ws.Range("D9").Value = '0.06723'
$ws.Range("E9").Value = '  +1.37%  '
$ws.Range("D10").Value = '1.898.69'
$ws.Range("E10").Value = '  +1.03%  '
$ws.Range("D11").Value = '17.21'
$ws.Range("E11").Value = '  +3.31%  '
$ws.Range("D12").Value = '0.07358'
$ws.Range("E12").Value = '  +1.77%  '
$ws.Range("D13").Value = '5.156'
$ws.Range("E13").Value = '  +3.99%  '
$ws.Range("D14").Value = '88.35'
$ws.Range("E14").Value = '  +0.20%  '
$ws.Range("D15").Value = '0.6704'
$ws.Range("E15").Value = '  +1.59%  '
$ws.Range("D16").Value = '30.561.19'
$ws.Range("E16").Value = '  +0.25%  '
$ws.Range("D17").Value = '0.000007882'
$ws.Range("E17").Value = '  +1.16%  '
$ws.Range("D18").Value = '13.46'
$ws.Range("E18").Value = '  +4.70%  '
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").Value = '2.159.39'
$ws.Range("E20").Value = '  +1.98%  '
$ws.Range("D21").Value = '5.322'
$ws.Range("E21").Value = '  +12.90%  '
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").Value = '191.80'
$ws.Range("E23").Value = '  +3.25%  '
$ws.Range("D24").Value = '6.228'
$ws.Range("E24").Value = '  +3.65%  '
$ws.Range("D25").Value = '9.560'
$ws.Range("E25").Value = '  +3.66%  '
$ws.Range("D26").Value = '160.66'
$ws.Range("E26").Value = '  +3.25%  '
$ws.Range("E27").Value = '  +0.69%  '
$ws.Range("D28").Value = '1.944'
$ws.Range("E28").Value = '  +6.32%  '
$ws.Range("E29").Value = '  +5.36%  '
$ws.Range("D30").Value = '4.458'
$ws.Range("E30").Value = '  +5.56%  '
$ws.Range("D31").Value = '0.09202'
$ws.Range("E31").Value = '  +2.42%  '
$ws.Range("D32").Value = '4.171'
$ws.Range("E32").Value = '  +6.92%  '
$ws.Range("D33").Value = '0.05241'
$ws.Range("E33").Value = '  +1.15%  '
$ws.Range("D34").Value = '0.7467'
$ws.Range("E34").Value = '  +2.27%  '
$ws.Range("D35").Value = '1.106'
$ws.Range("E35").Value = '  +3.17%  '
$ws.Range("D36").Value = '2.713'
$ws.Range("E36").Value = '  +0.50%  '
$ws.Range("D37").Value = '0.01836'
$ws.Range("E37").Value = '  +1.50%  '
$ws.Range("D38").Value = '2.697'
$ws.Range("E38").Value = '  +1.98%  '
$ws.Range("D39").Value = '0.9233'
$ws.Range("E39").Value = '  +0.68%  '
$ws.Range("D40").Value = '2.067'
$ws.Range("E40").Value = '  +2.02%  '
$ws.Range("D41").Value = '0.4434'
$ws.Range("E41").Value = '  +3.24%  '
$ws.Range("D42").Value = '5.960'
$ws.Range("E42").Value = '  +5.06%  '
$ws.Range("D43").Value = '71.99'
$ws.Range("E43").Value = '  +27.08%  '
$ws.Range("D44").Value = '106.29'
$ws.Range("E44").Value = '  +2.88%  '
$ws.Range("D45").Value = '0.9938'
$ws.Range("D46").Value = '0.1387'
$ws.Range("E46").Value = '  +4.41%  '
$ws.Range("D47").Value = '7.597'
$ws.Range("E47").Value = '  +5.41%  '
$ws.Range("D48").Value = '9.061'
$ws.Range("E48").Value = '  +6.09%  '
$ws.Range("D49").Value = '35.12'
$ws.Range("E49").Value = '  +6.28%  '
$ws.Range("D50").Value = '0.05832'
$ws.Range("E50").Value = '  +0.44%  '
$ws.Range("D51").Value = '0.3966'
$ws.Range("E51").Value = '  +2.59%  '

# Restore the original (unstyled) appearance of the edited range.
$ws.Range("D2:E51").Style = "Normal"
